# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (plain decimal strings like "297.93") need a transient text format so the
# literal string -- not a rounded/binary-float number -- is what gets stored;
# the format is reverted to Normal/General right after so the cell style is
# left exactly as it was (unstyled), matching the source data.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '39.831.68'
$ws.Range("E2").Value = '  -0.08%  '

$ws.Range("D3").Value = '2.224.73'
$ws.Range("E3").Value = '  -4.69%  '

$ws.Range("E4").Value = '  +0.09%  '

Set-TextValue $ws.Range("D5") '297.93'
$ws.Range("E5").Value = '  -3.31%  '

Set-TextValue $ws.Range("D6") '83.93'
$ws.Range("E6").Value = '  -0.49%  '

Set-TextValue $ws.Range("D7") '0.513'
$ws.Range("E7").Value = '  -2.47%  '

$ws.Range("E8").Value = '  +0.04%  '

Set-TextValue $ws.Range("D9") '0.467'
$ws.Range("E9").Value = '  -2.95%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range("D10") '0.0779'
$ws.Range("E10").Value = '  -4.76%  '

$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range("D11") '29.81'
$ws.Range("E11").Value = '  +0.22%  '

Set-TextValue $ws.Range("D12") '46.75'
$ws.Range("E12").Value = '  -10.46%  '

$ws.Range("E13").Value = '  -2.10%  '

$ws.Range("D14").Value = '2.572.18'
$ws.Range("E14").Value = '  -4.43%  '

Set-TextValue $ws.Range("D15") '6.30'
$ws.Range("E15").Value = '  -1.31%  '

Set-TextValue $ws.Range("D16") '14.09'
$ws.Range("E16").Value = '  -4.18%  '

$ws.Range("D17").Value = '2.232.66'
$ws.Range("E17").Value = '  -4.62%  '

Set-TextValue $ws.Range("D18") '0.717'
$ws.Range("E18").Value = '  -4.87%  '

$ws.Range("D19").Value = '39.762.59'
$ws.Range("E19").Value = '  -0.17%  '

$ws.Range("D20").Value = '0.0₃0879'
$ws.Range("E20").Value = '  -2.04%  '

Set-TextValue $ws.Range("D21") '5.79'
$ws.Range("E21").Value = '  -5.01%  '

Set-TextValue $ws.Range("D22") '10.53'
$ws.Range("E22").Value = '  +0.15%  '

Set-TextValue $ws.Range("D23") '65.06'
$ws.Range("E23").Value = '  -3.93%  '

Set-TextValue $ws.Range("D24") '233.55'
$ws.Range("E24").Value = '  -0.79%  '

$ws.Range("E25").Value = '  -0.16%  '

Set-TextValue $ws.Range("D26") '2.42'
$ws.Range("E26").Value = '  -4.52%  '

Set-TextValue $ws.Range("D27") '1.80'
$ws.Range("E27").Value = '  +0.73%  '

Set-TextValue $ws.Range("D28") '22.75'
$ws.Range("E28").Value = '  -1.88%  '

$ws.Range("E29").Value = '  +3.17%  '

Set-TextValue $ws.Range("D30") '9.17'
$ws.Range("E30").Value = '  -0.13%  '

Set-TextValue $ws.Range("D31") '32.51'
$ws.Range("E31").Value = '  -2.56%  '

Set-TextValue $ws.Range("D32") '149.34'
$ws.Range("E32").Value = '  -1.94%  '

$ws.Range("E33").Value = '  -0.12%  '

Set-TextValue $ws.Range("D34") '4.82'
$ws.Range("E34").Value = '  -4.75%  '

Set-TextValue $ws.Range("D35") '2.40'
$ws.Range("E35").Value = '  -1.57%  '

Set-TextValue $ws.Range("D36") '0.0700'
$ws.Range("E36").Value = '  -2.00%  '

Set-TextValue $ws.Range("D37") '16.30'
$ws.Range("E37").Value = '  +6.02%  '

Set-TextValue $ws.Range("D38") '0.111'
$ws.Range("E38").Value = '  -2.39%  '

Set-TextValue $ws.Range("D39") '0.0979'
$ws.Range("E39").Value = '  -0.93%  '

Set-TextValue $ws.Range("D40") '2.67'
$ws.Range("E40").Value = '  -2.19%  '

Set-TextValue $ws.Range("D41") '1.66'
$ws.Range("E41").Value = '  -3.03%  '

Set-TextValue $ws.Range("D42") '3.68'
$ws.Range("E42").Value = '  -2.70%  '

$ws.Range("D43").Value = '1.928.08'
$ws.Range("E43").Value = '  -0.40%  '

$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range("D44") '2.15'
$ws.Range("E44").Value = '  -2.90%  '

$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D45") '0.0266'
$ws.Range("E45").Value = '  +1.65%  '

$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D46") '9.28'
$ws.Range("E46").Value = '  +0.26%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D47") '16.55'
$ws.Range("E47").Value = '  -5.37%  '

Set-TextValue $ws.Range("D48") '2.61'
$ws.Range("E48").Value = '  -2.39%  '

$ws.Range("D49").Value = '2.442.08'
$ws.Range("E49").Value = '  -4.24%  '

Set-TextValue $ws.Range("D50") '70.95'
$ws.Range("E50").Value = '  +1.17%  '

Set-TextValue $ws.Range("D51") '88.88'
$ws.Range("E51").Value = '  -3.53%  '
